# "Add design type field to categories"
#
# 1. Categories sheet: add a "Design Type" column (B) with a value for every
#    existing category row, two brand-new category rows (International Union
#    Logos / Inspiration Board), and a list-data-validation on column B bound
#    to the DesignType named range.
# 2. Subcategories sheet: add the child rows for the two new categories
#    (the International Union Logos unions + the Inspiration Board groups),
#    extend the "Parent > Child" helper formula down through the new rows,
#    and tighten the existing list validation back down to a single
#    contiguous block.
# 3. Leave the workbook positioned on the Subcategories sheet/cell the
#    author ended the session on.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Categories sheet
# ---------------------------------------------------------------------------
$categories = $wb.Worksheets.Item("Categories")

$categories.Range("B1").Value = "Design Type"
$categories.Range("B2").Value = "Screen Print"
$categories.Range("B3").Value = "Screen Print"
$categories.Range("B4").Value = "Screen Print"
$categories.Range("B5").Value = "Screen Print"
$categories.Range("B6").Value = "Screen Print"

$categories.Range("A7").Value = "International Union Logos"
$categories.Range("B7").Value = "Embroidery"
$categories.Range("A8").Value = "Inspiration Board"
$categories.Range("B8").Value = "Embroidery"

$categories.Range("B2:B1048576").Validation.Add(3, 1, 1, "=DesignType")

# ---------------------------------------------------------------------------
# Subcategories sheet
# ---------------------------------------------------------------------------
$subcategories = $wb.Worksheets.Item("Subcategories")

$unions = @(
    "APWU",
    "BAC",
    "Boilermakers",
    "Carpenters",
    "CWA",
    "Electrical Workers",
    "IAFF",
    "IAM",
    "IATSE",
    "Ironworkers",
    "Heat & Frost Insulators",
    "IUEC",
    "IUOE",
    "IUPAT",
    "Laborers",
    "Mailhandlers",
    "NATCA",
    "NALC",
    "NEA",
    "OPCMIA",
    "OPEIU",
    "PASS",
    "Nurses",
    "Roofers",
    "SEIU",
    "SMART",
    "Teamsters",
    "UA",
    "UAW",
    "UWUA",
    "USW"
)

$row = 53
foreach ($union in $unions) {
    $subcategories.Range("A$row").Value = $union
    $subcategories.Range("B$row").Value = "International Union Logos"
    $row = $row + 1
}

$inspiration = @(
    "Hats/Beanies",
    "Chest/Sleeve",
    "Full Size Embroidery",
    "Mixed Media/Dye Sublimation/Applique"
)

foreach ($item in $inspiration) {
    $subcategories.Range("A$row").Value = $item
    $subcategories.Range("B$row").Value = "Inspiration Board"
    $row = $row + 1
}

# Extend the shared "Parent > Child" helper formula down through row 87.
$subcategories.Range("C53:C87").Formula = '=_xlfn.CONCAT(B53, " > ", A53)'

# Re-collapse the list validation on column B back to one contiguous block.
$subcategories.Range("B2:B1048576").Validation.Delete()
$subcategories.Range("B2:B1048576").Validation.Add(3, 1, 1, "=Categories")

# Column widths follow the new, wider content.
$categories.Columns.Item(1).ColumnWidth = 26
$subcategories.Columns.Item(1).ColumnWidth = 38
$subcategories.Columns.Item(3).ColumnWidth = 55

# ---------------------------------------------------------------------------
# Leave the workbook on the Subcategories sheet, scrolled to where the new
# rows were entered.
# ---------------------------------------------------------------------------
$designs = $wb.Worksheets.Item("Designs")
$designs.Activate()
$designs.Range("E2").Select()

$subcategories.Activate()
$subcategories.Range("F57").Select()
